# Update "想去人数" (want-to-go count) values in column F across the
# four worksheets, reflecting refreshed counts from the data source.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 885
$ws1.Range("F3").Value = 1467
$ws1.Range("F4").Value = 1123
$ws1.Range("F5").Value = 532
$ws1.Range("F8").Value = 685
$ws1.Range("F9").Value = 260
$ws1.Range("F11").Value = 99
$ws1.Range("F12").Value = 224
$ws1.Range("F13").Value = 162
$ws1.Range("F14").Value = 3059
$ws1.Range("F16").Value = 12
$ws1.Range("F17").Value = 440
$ws1.Range("F19").Value = 510
$ws1.Range("F20").Value = 285
$ws1.Range("F28").Value = 55
$ws1.Range("F29").Value = 1599
$ws1.Range("F30").Value = 331

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 240
$ws2.Range("F7").Value = 234
$ws2.Range("F12").Value = 134

# Sheet: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 73

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 885
$ws4.Range("F4").Value = 1467
$ws4.Range("F5").Value = 1123
$ws4.Range("F8").Value = 73
$ws4.Range("F9").Value = 532
$ws4.Range("F12").Value = 685
$ws4.Range("F14").Value = 260
$ws4.Range("F16").Value = 99
$ws4.Range("F17").Value = 224
$ws4.Range("F18").Value = 162
$ws4.Range("F19").Value = 3059
$ws4.Range("F21").Value = 12
$ws4.Range("F22").Value = 240
$ws4.Range("F23").Value = 440
$ws4.Range("F25").Value = 510
$ws4.Range("F26").Value = 285
$ws4.Range("F31").Value = 234
$ws4.Range("F37").Value = 134
$ws4.Range("F41").Value = 55
$ws4.Range("F42").Value = 1599
$ws4.Range("F43").Value = 331

$wb.Save()
